$wb = $excel.ActiveWorkbook

# --- Sheet 1 ("Metadata"): update Date and FHIR Version values ---
$wsMeta = $wb.Worksheets.Item(1)
$wsMeta.Range("B8").Value() = "2025-06-13T15:45:04+00:00"
$wsMeta.Range("B15").Value() = "4.0.1"

# --- Sheet 2 ("Elements"): update constraint text, type text, and URL text ---
$wsElem = $wb.Worksheets.Item(2)

$newConstraint = "ele-1:All FHIR elements must have a @value or children {hasValue() or (children().count() > id.count())}`next-1:Must have either extensions or value[x], not both {extension.exists() != value.exists()}"

# Row for "Extension" (row 2) constraint text loses the "unless an empty Parameters..." wording
$wsElem.Range("AJ2").Value() = $newConstraint

# Row for "Extension.extension" (row 4) already had this text; re-set so it now
# collapses (dedupes) onto the same shared string as AJ2, removing the old duplicate entry
$wsElem.Range("AJ4").Value() = $newConstraint

# Row for "Extension.id" (row 3): Type(s) changes from "id" to "string"
$wsElem.Range("K3").Value() = "string`n"

# Row for "Extension.value[x]" (row 6): Definition URL changes from R4B to R4
$wsElem.Range("M6").Value() = "Value of extension - must be one of a constrained set of the data types (see [Extensibility](http://hl7.org/fhir/R4/extensibility.html) for a list)."
